$wb = $excel.ActiveWorkbook

# --- Summary sheet: update "Generated" timestamp ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B2").Value = "2026-02-04 18:00:32"

# --- Monthly_EPU sheet: fill in ireland_domestic_epu for Jun-Oct 2025 ---
$monthly = $wb.Worksheets.Item("Monthly_EPU")
$monthly.Range("F2").Value = 57.36575675314314
$monthly.Range("F3").Value = 94.69138610935934
$monthly.Range("F4").Value = -54.63860549041632
$monthly.Range("F5").Value = -16.55135733179955
$monthly.Range("F6").Value = 35.12546384883223

# --- Quarterly_EPU sheet: fill in ireland_domestic_epu for 2025Q2-Q4 ---
$quarterly = $wb.Worksheets.Item("Quarterly_EPU")
$quarterly.Range("C2").Value = 57.36575675314314
$quarterly.Range("C3").Value = 7.833807762381156
$quarterly.Range("C4").Value = 35.12546384883223
